$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column S: header "Es apto" with same formatting as header row (R1) ---
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S1").Value = "Es apto"

# --- Column S rows 2-5: value "S" with same formatting as data rows (R2) ---
$ws.Range("R2").Copy()
$ws.Range("S2:S5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S2").Value = "S"
$ws.Range("S3").Value = "S"
$ws.Range("S4").Value = "S"
$ws.Range("S5").Value = "S"

# --- Column S row 6: value "N" with same formatting as R6 ---
$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S6").Value = "N"

$excel.CutCopyMode = 0

# --- Row 7: new unstyled data row (template header row for html) ---
$ws.Range("D7").Value = "IB"
$ws.Range("E7").Value = "Interno"
$ws.Range("F7").Value = "High Reach"
$ws.Range("J7").Value = "S"
$ws.Range("K7").Value = "S"
$ws.Range("L7").Value = "S"
$ws.Range("N7").Value = "S"
$ws.Range("O7").Value = " "
$ws.Range("P7").Value = "S"
$ws.Range("Q7").Value = "S"
$ws.Range("R7").Value = "S"

# --- Update selection to match the new used range ---
$ws.Range("A2:S6").Select()
